# issue #5: add legislator_id, name, date into dataframe
#
# Adds three columns - date, legislator_name, legislator_id - to the
# "股票" (stock) sheet (the 4th worksheet) of the property-declaration
# workbook, matching the per-legislator metadata encoded in the source
# filename: 林明溱 / 2012-04-18 / id 1706.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

$legislatorDate = "2012-04-18"
$legislatorName = "林明溱"
$legislatorId   = 1706

# --- Header row (H1:J1) -----------------------------------------------
# Copy the existing header formatting (bold / bordered / centered) from
# G1 onto the new header cells *before* writing their text, so the COM
# bridge resolves them to the already-existing header style instead of
# minting new ones.
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("H1").Value = "date"
$ws.Range("I1").Value = "legislator_name"
$ws.Range("J1").Value = "legislator_id"

# --- Data rows -----------------------------------------------------
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    # The date column holds a literal text value ("2012-04-18"), not a
    # real date serial. Route it through a text formula first and then
    # collapse the formula to its literal string result, so Excel's
    # "looks like a date" auto-conversion never kicks in.
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = '="' + $legislatorDate + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null

    $ws.Cells.Item($r, 9).Value  = $legislatorName
    $ws.Cells.Item($r, 10).Value = $legislatorId
}
